$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells we are about to rewrite to keep their literal
# text representation (e.g. "73.80", "42.642.62") instead of being auto-parsed
# as numbers/dates by Excel. Apply as one unioned range so only a single shared
# style entry is introduced.
$dRange = $ws.Range("D2,D3,D5,D6,D7,D9,D10,D11,D12,D13,D15,D16,D17,D18,D19,D21,D22,D23,D24,D26,D27,D28,D29,D31,D32,D33,D35,D36,D38,D40,D41,D43,D45,D46,D47,D51")
foreach ($area in $dRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = '42.642.62'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '2.283.11'
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '251.86'
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").Value = '0.637'
$ws.Range("E6").Value = '  +2.47%  '
$ws.Range("D7").Value = '73.73'
$ws.Range("E7").Value = '  +9.29%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.647'
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("D10").Value = '39.15'
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("D11").Value = '0.0992'
$ws.Range("E11").Value = '  +5.98%  '
$ws.Range("D12").Value = '59.02'
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").Value = '7.35'
$ws.Range("E13").Value = '  +4.56%  '
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("D15").Value = '2.625.56'
$ws.Range("E15").Value = '  +3.44%  '
$ws.Range("D16").Value = '14.98'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").Value = '0.874'
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").Value = '2.276.62'
$ws.Range("E18").Value = '  +3.25%  '
$ws.Range("D19").Value = '42.576.82'
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("E20").Value = '  +5.15%  '
$ws.Range("D21").Value = '6.31'
$ws.Range("E21").Value = '  +2.85%  '
$ws.Range("D22").Value = '72.33'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '234.99'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  +9.87%  '
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '11.48'
$ws.Range("E26").Value = '  +3.11%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.20%  '
$ws.Range("D28").Value = '2.42'
$ws.Range("E28").Value = '  +0.83%  '
$ws.Range("D29").Value = '3.66'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '166.90'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '21.08'
$ws.Range("E32").Value = '  +3.44%  '
$ws.Range("D33").Value = '6.37'
$ws.Range("E33").Value = '  +8.39%  '
$ws.Range("E34").Value = '  +5.65%  '
$ws.Range("D35").Value = '0.0819'
$ws.Range("E35").Value = '  +5.19%  '
$ws.Range("D36").Value = '32.23'
$ws.Range("E36").Value = '  +24.02%  '
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("D38").Value = '4.73'
$ws.Range("E38").Value = '  +16.22%  '
$ws.Range("E39").Value = '  +3.63%  '
$ws.Range("D40").Value = '0.0307'
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").Value = '13.85'
$ws.Range("E41").Value = '  +15.83%  '
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("D43").Value = '5.96'
$ws.Range("E43").Value = '  +6.03%  '
$ws.Range("E44").Value = '  +9.01%  '
$ws.Range("D45").Value = '9.15'
$ws.Range("E45").Value = '  +6.96%  '
$ws.Range("D46").Value = '61.96'
$ws.Range("D47").Value = '4.83'
$ws.Range("E47").Value = '  -6.48%  '
$ws.Range("E48").Value = '  +4.41%  '
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("D51").Value = '97.46'
$ws.Range("E51").Value = '  +4.91%  '
